$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 12:20"

# Row 10
$ws.Range("B10").Value = 50468
$ws.Range("C10").Value = 2875
$ws.Range("D10").Value = 16711
$ws.Range("E10").Value = 30597
$ws.Range("F10").Value = 3956
$ws.Range("G10").Value = 124
$ws.Range("H10").Value = 3160

# Row 12
$ws.Range("B12").Value = 18117
$ws.Range("C12").Value = 349
$ws.Range("D12").Value = 4013
$ws.Range("E12").Value = 13599
$ws.Range("F12").Value = 348
$ws.Range("G12").Value = 17
$ws.Range("H12").Value = 505

# Row 16
$ws.Range("B16").Value = 10877
$ws.Range("C16").Value = 166
$ws.Range("D16").Value = 1749
$ws.Range("E16").Value = 8970
$ws.Range("F16").Value = 227
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 158

# Row 32 -> Rumania
$ws.Range("A32").Value = "Rumania"
$ws.Range("B32").Value = 2738
$ws.Range("C32").Value = 278
$ws.Range("D32").Value = 267
$ws.Range("E32").Value = 2377
$ws.Range("F32").Value = 78
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = 94

# Row 33 -> Polonia
$ws.Range("A33").Value = "Polonia"
$ws.Range("B33").Value = 2633
$ws.Range("C33").Value = 79
$ws.Range("D33").Value = 56
$ws.Range("E33").Value = 2532
$ws.Range("F33").Value = 50
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 45

# Row 34 -> Filipinas
$ws.Range("A34").Value = "Filipinas"
$ws.Range("B34").Value = 2633
$ws.Range("C34").Value = 322
$ws.Range("D34").Value = 51
$ws.Range("E34").Value = 2475
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 11
$ws.Range("H34").Value = 107

# Row 61
$ws.Range("B61").Value = 802
$ws.Range("C61").Value = 36
$ws.Range("D61").Value = 154
$ws.Range("E61").Value = 644
$ws.Range("F61").Value = 8
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 4

# Row 75 -> Eslovaquia
$ws.Range("A75").Value = "Eslovaquia"
$ws.Range("B75").Value = 426
$ws.Range("C75").Value = 26
$ws.Range("D75").Value = 3
$ws.Range("E75").Value = 422
$ws.Range("F75").Value = 1
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 1

# Row 76 -> Tunez
$ws.Range("A76").Value = "Tunez"
$ws.Range("B76").Value = 423
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 5
$ws.Range("E76").Value = 406
$ws.Range("F76").Value = 10
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 12

# Row 77 -> Moldavia
$ws.Range("A77").Value = "Moldavia"
$ws.Range("B77").Value = 423
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 23
$ws.Range("E77").Value = 395
$ws.Range("F77").Value = 65
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 5

# Row 78 -> Kazajistan
$ws.Range("A78").Value = "Kazajistan"
$ws.Range("B78").Value = 402
$ws.Range("C78").Value = 22
$ws.Range("D78").Value = 27
$ws.Range("E78").Value = 372
$ws.Range("F78").Value = 6
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 3

# Row 99 -> Senegal
$ws.Range("A99").Value = "Senegal"
$ws.Range("B99").Value = 195
$ws.Range("C99").Value = 5
$ws.Range("D99").Value = 55
$ws.Range("E99").Value = 139
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 1

# Row 100 -> Costa de Marfil
$ws.Range("A100").Value = "Costa de Marfil"
$ws.Range("B100").Value = 190
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 9
$ws.Range("E100").Value = 180
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 1

# Row 101 -> Uzbekistan
$ws.Range("A101").Value = "Uzbekistan"
$ws.Range("B101").Value = 190
$ws.Range("C101").Value = 9
$ws.Range("D101").Value = 12
$ws.Range("E101").Value = 176
$ws.Range("F101").Value = 8
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 2

# Row 148
$ws.Range("B148").Value = 29
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 3
$ws.Range("E148").Value = 26
$ws.Range("F148").Value = 2
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 0

# Row 157 -> Nueva Caledonia
$ws.Range("A157").Value = "Nueva Caledonia"
$ws.Range("B157").Value = 18
$ws.Range("C157").Value = 2
$ws.Range("D157").Value = 1
$ws.Range("E157").Value = 17
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 0

# Row 158 -> Gabon
$ws.Range("A158").Value = "Gabon"
$ws.Range("B158").Value = 18
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 0
$ws.Range("E158").Value = 17
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 1

# Row 159 -> Islas Virgenes de los Estados Unidos
$ws.Range("A159").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B159").Value = 17
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 0
$ws.Range("E159").Value = 17
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 0
